$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 18: sv/Statement-opinion -> sd/Statement-non-opinion
$ws.Range("I18").Value = "sd"
$ws.Range("J18").Value = "Statement-non-opinion"

# Row 31: %/Uninterpretable -> sd/Statement-non-opinion
$ws.Range("I31").Value = "sd"
$ws.Range("J31").Value = "Statement-non-opinion"

# Row 32: %/Uninterpretable -> sd/Statement-non-opinion
$ws.Range("I32").Value = "sd"
$ws.Range("J32").Value = "Statement-non-opinion"

# Row 35: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I35").Value = "sv"
$ws.Range("J35").Value = "Statement-opinion"

# Row 40: sv/Statement-opinion -> sd/Statement-non-opinion
$ws.Range("I40").Value = "sd"
$ws.Range("J40").Value = "Statement-non-opinion"

# Row 43: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I43").Value = "sv"
$ws.Range("J43").Value = "Statement-opinion"

# Row 58: %/Uninterpretable -> sd/Statement-non-opinion
$ws.Range("I58").Value = "sd"
$ws.Range("J58").Value = "Statement-non-opinion"

# Row 59: %/Uninterpretable -> sd/Statement-non-opinion
$ws.Range("I59").Value = "sd"
$ws.Range("J59").Value = "Statement-non-opinion"
